# Update asset numbers (ANLN1) and asset sub-numbers (ANLN2) on the "Data"
# sheet of the MPA test automation upload file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Column D (ANLN1) updates: 60000468 -> 60000476 for the group of rows
# sharing the same original asset number.
$ws.Range("D6").Value  = 60000476
$ws.Range("D7").Value  = 60000476
$ws.Range("D8").Value  = 60000476
$ws.Range("D9").Value  = 60000476
$ws.Range("D10").Value = 60000476
$ws.Range("D16").Value = 60000476
$ws.Range("D17").Value = 60000476

$ws.Range("D20").Value = 60000477
$ws.Range("D22").Value = 60000478
$ws.Range("D24").Value = 60000479
$ws.Range("D26").Value = 60000480

# Column E (ANLN2) updates: shift sub-number values by +7.
$ws.Range("E11").Value = 326
$ws.Range("E12").Value = 326
$ws.Range("E13").Value = 326
$ws.Range("E14").Value = 326
$ws.Range("E15").Value = 326
$ws.Range("E18").Value = 326
$ws.Range("E19").Value = 326

$ws.Range("E21").Value = 327
$ws.Range("E23").Value = 328
$ws.Range("E25").Value = 329
$ws.Range("E27").Value = 330
